$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.840.25"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = "'1.768.32"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('D4').Value = "'1.004"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').Value = "'327.72"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').Value = "'1.003"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('D7').Value = "'0.4483"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.81%  '
$ws.Range('D8').Value = "'0.3566"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.14%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'42.12"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.07426"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.51%  '
$ws.Range('D11').Value = "'1.093"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = "'20.84"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = "'6.038"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').Value = "'7.209"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').Value = "'1.775.79"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range('D17').Value = "'92.67"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = "'0.00001058"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('D19').Value = "'0.06409"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = "'17.17"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('D22').Value = "'5.813"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = "'27.870.34"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('D24').Value = "'11.30"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').Value = "'2.117"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').Value = "'162.43"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').Value = "'20.21"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.50%  '
$ws.Range('D28').Value = "'1.976.63"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.07%  '
$ws.Range('D29').Value = "'2.182"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.72%  '
$ws.Range('D30').Value = "'125.37"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('D31').Value = "'1.098"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('D32').Value = "'0.09130"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').Value = "'5.571"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('D34').Value = "'3.636"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('D35').Value = "'11.78"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('D36').Value = "'0.02288"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').Value = "'0.06093"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').Value = "'0.2095"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').Value = "'0.6320"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').Value = "'4.950"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.73%  '
$ws.Range('D41').Value = "'1.183"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('D42').Value = "'1.394"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('D43').Value = "'7.908"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('D44').Value = "'13.20"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('D45').Value = "'3.740"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').Value = "'0.5857"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('D47').Value = "'122.42"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').Value = "'1.954"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.61%  '
$ws.Range('D49').Value = "'0.06907"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('D50').Value = "'1.137"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('D51').Value = "'72.82"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.53%  '
